$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Column D: "blur" second calc column ---
$ws.Range("D2").Value = "blur"

$dValues = @(45,58,39,104,37,24,199,176,181,138,93,208,149,2,30,2,156,86,92,106,71,76,11,22,191,8,208,196,188,190)
for ($i = 0; $i -lt $dValues.Count; $i++) {
    $row = 3 + $i
    $ws.Cells.Item($row, 4).Value = $dValues[$i]
}

# --- Column G: list of ids measured >= 145 ---
$ws.Range("G2").Value = "145 이상 측정"

$gValues = @("t07","t08","t09","t12","t13","t17","t25","t27","t28","t29","t30")
for ($i = 0; $i -lt $gValues.Count; $i++) {
    $row = 4 + $i
    $ws.Cells.Item($row, 7).Value = $gValues[$i]
}

# --- Column widths ---
# ColumnWidth is in characters and gets pixel-snapped internally; 13.29
# lands on the same pixel bucket as a stored width of 14.
$ws.Columns.Item(7).ColumnWidth = 13.29

# --- View settings ---
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B15").Select()

# Workbook window position (best effort; matches the author's on-screen move)
$wb.Windows.Item(1).Left = 6300
$wb.Windows.Item(1).Top = 3420
